# Modification du journal et création pour le désigne des pages web
#
# 1) "Le logo du site" -> append " avec le logo du Fairmont le Montreux
#    Palace" as its own run, and move the _GoBack bookmark here (it used to
#    sit at the very end of the document).
# 2) Drop stale w:proofErr spellStart/spellEnd wrappers around a few
#    single-word runs ("N°Ticket", "Capchat", the login page's "Username"
#    and "Password").
# 3) Collapse the three runs (with a gramStart/gramEnd proofErr pair in the
#    middle) that spell out "Nom du machine qui est en panne" into one run.

$d = $word.ActiveDocument

function Remove-ProofErrAroundParagraph($paragraphIndex) {
    $p = $d.Paragraphs($paragraphIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $text = $rng.Text
    $rng.Delete()
    $rng.InsertAfter($text)
}

# --- "N°Ticket" (paragraph 7) ---------------------------------------------
Remove-ProofErrAroundParagraph 7

# --- "Nom du machine qui est en panne" (paragraph 27) ----------------------
# Originally split across 3 runs ("Nom " / "du" / " machine qui est en
# panne") with a gramStart/gramEnd proofErr pair wrapping the middle run.
# Re-write it as a single clean run.
$p27 = $d.Paragraphs(27)
$rng27 = $p27.Range
$rng27.MoveEnd(1, -1) | Out-Null
$rng27.Delete()
$rng27.InsertAfter("Nom du machine qui est en panne")

# --- "Capchat" (paragraph 32) ----------------------------------------------
Remove-ProofErrAroundParagraph 32

# --- "Username" / "Password" on the Login page (paragraphs 41 & 42) -------
Remove-ProofErrAroundParagraph 41
Remove-ProofErrAroundParagraph 42

# --- "Le logo du site" (paragraph 3): add the extra run + move bookmark ---
$p3 = $d.Paragraphs(3)
$rng3 = $p3.Range
$rng3.MoveEnd(1, -1) | Out-Null
$endOfText = $rng3.End

# Split the paragraph so the new text lands in its own run (rather than
# being merged into the existing "Le logo du site" run), then stitch the
# paragraph back together by deleting the paragraph mark that separated
# them.
$splitPoint = $d.Range($endOfText, $endOfText)
$splitPoint.InsertParagraphAfter()
$newPara = $p3.Next()
$newPara.Range.InsertAfter(" avec le logo du Fairmont le Montreux Palace")
$pilcrow = $d.Range($p3.Range.End - 1, $p3.Range.End)
$pilcrow.Delete()

# Relocate the _GoBack bookmark to sit right after the text we just added.
# (Word keeps exactly one bookmark per name, so adding it here also removes
# it from its old location at the end of the document; we also delete it
# explicitly first to be safe.)
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$rng3b = $p3.Range
$rng3b.MoveEnd(1, -1) | Out-Null
$paraEnd = $rng3b.End
$marker = $d.Range($paraEnd, $paraEnd)
$marker.InsertAfter("@")
$bookmarkRange = $d.Range($paraEnd, $paraEnd + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
$d.Range($paraEnd, $paraEnd + 1).Delete()
